$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the 9 sending/target cluster combinations (ECs/FAPs/sCs) x (ECs/FAPs/sCs)
# for the Dcn -> Tlr4 ligand-receptor pair, replacing the previous 1x1x1 block of rows 2-7
# with the corrected 3x3 block spanning rows 2-10, per updated NATMI computation.

# Row 2: ECs -> ECs
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Dcn"
$ws.Cells.Item(2, 3).Value = "Tlr4"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 6.615074666666668
$ws.Cells.Item(2, 8).Value = 19.845224
$ws.Cells.Item(2, 9).Value = 0.0008916467884469992
$ws.Cells.Item(2, 10).Value = 0.0008916467884469989
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 22.48784766666667
$ws.Cells.Item(2, 14).Value = 67.463543
$ws.Cells.Item(2, 15).Value = 0.4520839499795984
$ws.Cells.Item(2, 16).Value = 0.4520839499795983
$ws.Cells.Item(2, 17).Value = 148.7587914076258
$ws.Cells.Item(2, 18).Value = 1338.829122668632
$ws.Cells.Item(2, 19).Value = 0.0004030992021077427
$ws.Cells.Item(2, 20).Value = 0.0004030992021077425

# Row 3: ECs -> FAPs
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Dcn"
$ws.Cells.Item(3, 3).Value = "Tlr4"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 6.615074666666668
$ws.Cells.Item(3, 8).Value = 19.845224
$ws.Cells.Item(3, 9).Value = 0.0008916467884469992
$ws.Cells.Item(3, 10).Value = 0.0008916467884469989
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 23.26810333333333
$ws.Cells.Item(3, 14).Value = 69.80431
$ws.Cells.Item(3, 15).Value = 0.4677698025791556
$ws.Cells.Item(3, 16).Value = 0.4677698025791556
$ws.Cells.Item(3, 17).Value = 153.9202409017156
$ws.Cells.Item(3, 18).Value = 1385.28216811544
$ws.Cells.Item(3, 19).Value = 0.0004170854422021909
$ws.Cells.Item(3, 20).Value = 0.0004170854422021908

# Row 4: ECs -> sCs
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Dcn"
$ws.Cells.Item(4, 3).Value = "Tlr4"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 6.615074666666668
$ws.Cells.Item(4, 8).Value = 19.845224
$ws.Cells.Item(4, 9).Value = 0.0008916467884469992
$ws.Cells.Item(4, 10).Value = 0.0008916467884469989
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 3.986685666666666
$ws.Cells.Item(4, 14).Value = 11.960057
$ws.Cells.Item(4, 15).Value = 0.08014624744124609
$ws.Cells.Item(4, 16).Value = 0.08014624744124607
$ws.Cells.Item(4, 17).Value = 26.37222335752978
$ws.Cells.Item(4, 18).Value = 237.350010217768
$ws.Cells.Item(4, 19).Value = 0.00007146214413706559
$ws.Cells.Item(4, 20).Value = 0.00007146214413706557

# Row 5: FAPs -> ECs
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Dcn"
$ws.Cells.Item(5, 3).Value = "Tlr4"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 7285.701009
$ws.Cells.Item(5, 8).Value = 21857.103027
$ws.Cells.Item(5, 9).Value = 0.9820406017477925
$ws.Cells.Item(5, 10).Value = 0.9820406017477923
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 22.48784766666667
$ws.Cells.Item(5, 14).Value = 67.463543
$ws.Cells.Item(5, 15).Value = 0.4520839499795984
$ws.Cells.Item(5, 16).Value = 0.4520839499795983
$ws.Cells.Item(5, 17).Value = 163839.7344352716
$ws.Cells.Item(5, 18).Value = 1474557.609917445
$ws.Cells.Item(5, 19).Value = 0.4439647942784837
$ws.Cells.Item(5, 20).Value = 0.4439647942784835

# Row 6: FAPs -> FAPs
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Dcn"
$ws.Cells.Item(6, 3).Value = "Tlr4"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 7285.701009
$ws.Cells.Item(6, 8).Value = 21857.103027
$ws.Cells.Item(6, 9).Value = 0.9820406017477925
$ws.Cells.Item(6, 10).Value = 0.9820406017477923
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 23.26810333333333
$ws.Cells.Item(6, 14).Value = 69.80431
$ws.Cells.Item(6, 15).Value = 0.4677698025791556
$ws.Cells.Item(6, 16).Value = 0.4677698025791556
$ws.Cells.Item(6, 17).Value = 169524.4439331829
$ws.Cells.Item(6, 18).Value = 1525719.995398646
$ws.Cells.Item(6, 19).Value = 0.4593689384042801
$ws.Cells.Item(6, 20).Value = 0.4593689384042799

# Row 7: FAPs -> sCs
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Dcn"
$ws.Cells.Item(7, 3).Value = "Tlr4"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 7285.701009
$ws.Cells.Item(7, 8).Value = 21857.103027
$ws.Cells.Item(7, 9).Value = 0.9820406017477925
$ws.Cells.Item(7, 10).Value = 0.9820406017477923
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 3.986685666666666
$ws.Cells.Item(7, 14).Value = 11.960057
$ws.Cells.Item(7, 15).Value = 0.08014624744124609
$ws.Cells.Item(7, 16).Value = 0.08014624744124607
$ws.Cells.Item(7, 17).Value = 29045.79978419917
$ws.Cells.Item(7, 18).Value = 261412.1980577925
$ws.Cells.Item(7, 19).Value = 0.07870686906502879
$ws.Cells.Item(7, 20).Value = 0.07870686906502875

# Row 8: sCs -> ECs
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Dcn"
$ws.Cells.Item(8, 3).Value = "Tlr4"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 126.6246363333333
$ws.Cells.Item(8, 8).Value = 379.873909
$ws.Cells.Item(8, 9).Value = 0.01706775146376063
$ws.Cells.Item(8, 10).Value = 0.01706775146376063
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 22.48784766666667
$ws.Cells.Item(8, 14).Value = 67.463543
$ws.Cells.Item(8, 15).Value = 0.4520839499795984
$ws.Cells.Item(8, 16).Value = 0.4520839499795983
$ws.Cells.Item(8, 17).Value = 2847.515532711066
$ws.Cells.Item(8, 18).Value = 25627.63979439959
$ws.Cells.Item(8, 19).Value = 0.007716056499006978
$ws.Cells.Item(8, 20).Value = 0.007716056499006974

# Row 9: sCs -> FAPs
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Dcn"
$ws.Cells.Item(9, 3).Value = "Tlr4"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 126.6246363333333
$ws.Cells.Item(9, 8).Value = 379.873909
$ws.Cells.Item(9, 9).Value = 0.01706775146376063
$ws.Cells.Item(9, 10).Value = 0.01706775146376063
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 23.26810333333333
$ws.Cells.Item(9, 14).Value = 69.80431
$ws.Cells.Item(9, 15).Value = 0.4677698025791556
$ws.Cells.Item(9, 16).Value = 0.4677698025791556
$ws.Cells.Item(9, 17).Value = 2946.315122749755
$ws.Cells.Item(9, 18).Value = 26516.83610474779
$ws.Cells.Item(9, 19).Value = 0.007983778732673405
$ws.Cells.Item(9, 20).Value = 0.007983778732673402

# Row 10: sCs -> sCs
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Dcn"
$ws.Cells.Item(10, 3).Value = "Tlr4"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 126.6246363333333
$ws.Cells.Item(10, 8).Value = 379.873909
$ws.Cells.Item(10, 9).Value = 0.01706775146376063
$ws.Cells.Item(10, 10).Value = 0.01706775146376063
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 3.986685666666666
$ws.Cells.Item(10, 14).Value = 11.960057
$ws.Cells.Item(10, 15).Value = 0.08014624744124609
$ws.Cells.Item(10, 16).Value = 0.08014624744124607
$ws.Cells.Item(10, 17).Value = 504.8126227169792
$ws.Cells.Item(10, 18).Value = 4543.313604452813
$ws.Cells.Item(10, 19).Value = 0.001367916232080249
$ws.Cells.Item(10, 20).Value = 0.001367916232080249

